$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.607.62'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '2.516.19'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.587'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.32%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.538'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.62'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.55%  '
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").Value = '2.894.84'
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.99%  '
$ws.Range("D16").Value = '2.486.03'
$ws.Range("E16").Value = '  -2.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.858'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D18").Value = '42.586.84'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.67%  '
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.00%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("E28").Value = '  +11.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.37'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.34'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.46%  '
$ws.Range("E38").Value = '  -1.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.30'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.23%  '
$ws.Range("E40").Value = '  +1.27%  '
$ws.Range("E41").Value = '  +0.47%  '
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("E43").Value = '  -3.55%  '
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '2.025.10'
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.89%  '
$ws.Range("D49").Value = '2.755.51'
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.191'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.68%  '
